# Added filtering options for the Component Analysis
# Trim the trailing quarter-error columns on each data row so that each
# row only keeps the "near-term" forecast horizons that are still in
# range for that quarter, producing the staircase layout used by the
# filtered Component Analysis export.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> first column (letter) to clear through column K (last used column).
$rowsToClear = [ordered]@{
    2  = "G"
    3  = "I"
    4  = "G"
    5  = "I"
    6  = "G"
    7  = "I"
    8  = "G"
    9  = "I"
    10 = "G"
    11 = "I"
    12 = "G"
    13 = "I"
    14 = "G"
    15 = "I"
    16 = "G"
    17 = "I"
    18 = "K"
    19 = "I"
    20 = "K"
    21 = "I"
    22 = "K"
    23 = "J"
    24 = "I"
    26 = "K"
    27 = "J"
    28 = "I"
    30 = "K"
    31 = "J"
    32 = "I"
    34 = "K"
    35 = "J"
    36 = "I"
    38 = "K"
    39 = "J"
    40 = "I"
    42 = "K"
    43 = "J"
    44 = "I"
}

foreach ($row in $rowsToClear.Keys) {
    $firstCol = $rowsToClear[$row]
    $lastCol = "K"
    if ($row -eq 44) { $lastCol = "J" }
    $ws.Range("$firstCol$row`:$lastCol$row").ClearContents()
}
